$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44292
$ws.Range("M2").Value = 25
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 16000
$ws.Range("S2").Value = 2286
$ws.Range("D3").Value = 44292
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 2143
$ws.Range("D4").Value = 44301
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("Q4").Value = '$/bandeja 7 kilos'
$ws.Range("S4").Value = 2000
$ws.Range("T4").Value = 7
$ws.Range("D5").Value = 44301
$ws.Range("L5").Value = 'Segunda'
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("R5").Value = 'Región Metropolitana'
$ws.Range("S5").Value = 1714
$ws.Range("D6").Value = 44980
$ws.Range("L6").Value = 'Primera'
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 16000
$ws.Range("R6").Value = 'Región Metropolitana'
$ws.Range("S6").Value = 2286
$ws.Range("D7").Value = 44980
$ws.Range("L7").Value = 'Segunda'
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("S7").Value = 1857
$ws.Range("D8").Value = 44322
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 45
$ws.Range("D9").Value = 44322
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 1143
$ws.Range("D10").Value = 44302
$ws.Range("L10").Value = 'Primera'
$ws.Range("M10").Value = 50
$ws.Range("D11").Value = 44302
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 30
$ws.Range("D12").Value = 44299
$ws.Range("L12").Value = 'Primera'
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("R12").Value = 'Provincia de Santiago'
$ws.Range("S12").Value = 2143
$ws.Range("D13").Value = 44299
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 75
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("R13").Value = 'Provincia de Santiago'
$ws.Range("S13").Value = 1714
$ws.Range("D14").Value = 44320
$ws.Range("L14").Value = 'Primera'
$ws.Range("M14").Value = 20
$ws.Range("D15").Value = 44320
$ws.Range("L15").Value = 'Segunda'
$ws.Range("M15").Value = 30
$ws.Range("N15").Value = 8000
$ws.Range("O15").Value = 8000
$ws.Range("P15").Value = 8000
$ws.Range("S15").Value = 1143
$ws.Range("L16").Value = 'Primera'
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("S16").Value = 2143
$ws.Range("D17").Value = 44300
$ws.Range("L17").Value = 'Segunda'
$ws.Range("N17").Value = 12000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 12000
$ws.Range("S17").Value = 1714
$ws.Range("D18").Value = 44971
$ws.Range("L18").Value = 'Primera'
$ws.Range("M18").Value = 25
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("Q18").Value = '$/bandeja 5 kilos'
$ws.Range("S18").Value = 3000
$ws.Range("T18").Value = 5
